$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets updated (new "Neutro" sending cluster, Calca->Calcr, target sCs)
$ws.Range("A2").Value = "Neutro"
$ws.Range("B2").Value = "Calca"
$ws.Range("C2").Value = "Calcr"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1472293333333334
$ws.Range("H2").Value = 0.441688
$ws.Range("I2").Value = 0.2923521718179941
$ws.Range("J2").Value = 0.2923521718179941
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.194209333333333
$ws.Range("N2").Value = 15.582628
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.7647399773404445
$ws.Range("R2").Value = 6.882659796064
$ws.Range("S2").Value = 0.2923521718179941
$ws.Range("T2").Value = 0.2923521718179941

# New row 3 is added (sCs sending cluster, Calca->Calcr, target sCs)
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Calca"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3563733333333333
$ws.Range("H3").Value = 1.06912
$ws.Range("I3").Value = 0.7076478281820058
$ws.Range("J3").Value = 0.7076478281820059
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.194209333333333
$ws.Range("N3").Value = 15.582628
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.851077694151111
$ws.Range("R3").Value = 16.65969924736
$ws.Range("S3").Value = 0.7076478281820058
$ws.Range("T3").Value = 0.7076478281820059
